# Scheduled runner update: refresh currentAveragePrice / Leve profit figures
# across the per-job "Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 184.875
$ws.Range("I4").Value = 104.5
$ws.Range("K4").Value = 104.5
$ws.Range("M4").Value = 9.5
$ws.Range("H11").Value = 283.77777
$ws.Range("I11").Value = 283.77777
$ws.Range("K11").Value = 283.77777
$ws.Range("M11").Value = -143.77777
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H28").Value = 126043.875
$ws.Range("I28").Value = 154784.77
$ws.Range("J28").Value = 1500
$ws.Range("K28").Value = 154784.77
$ws.Range("L28").Value = 1500
$ws.Range("M28").Value = -154299.77
$ws.Range("N28").Value = -2470
$ws.Range("H111").Value = 1970.3572
$ws.Range("J111").Value = 2226.7144
$ws.Range("L111").Value = 6680.1432
$ws.Range("N111").Value = -12814.1432
$ws.Range("H116").Value = 23822384
$ws.Range("I116").Value = 31264794
$ws.Range("K116").Value = 31264794
$ws.Range("M116").Value = -31261352
$ws.Range("H132").Value = 1865
$ws.Range("I132").Value = 1873.5652
$ws.Range("J132").Value = 1799.3334
$ws.Range("K132").Value = 5620.6956
$ws.Range("L132").Value = 5398.0002
$ws.Range("M132").Value = -3090.6956
$ws.Range("N132").Value = -10458.0002
$ws.Range("H137").Value = 1356711.4
$ws.Range("I137").Value = 1369.4828
$ws.Range("J137").Value = 4632121
$ws.Range("K137").Value = 4108.4484
$ws.Range("L137").Value = 13896363
$ws.Range("M137").Value = -1558.4484
$ws.Range("N137").Value = -13901463
$ws.Range("H138").Value = 2038.19
$ws.Range("J138").Value = 2741.4614
$ws.Range("L138").Value = 8224.3842
$ws.Range("N138").Value = -18504.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3325.0715
$ws.Range("I2").Value = 3006.625
$ws.Range("J2").Value = 3749.6667
$ws.Range("K2").Value = 3006.625
$ws.Range("L2").Value = 3749.6667
$ws.Range("M2").Value = -2893.625
$ws.Range("N2").Value = -3975.6667
$ws.Range("H45").Value = 3254.3635
$ws.Range("I45").Value = 1249.5
$ws.Range("J45").Value = 3699.889
$ws.Range("K45").Value = 1249.5
$ws.Range("L45").Value = 3699.889
$ws.Range("M45").Value = -872.5
$ws.Range("N45").Value = -4453.889
$ws.Range("H61").Value = 2764.535
$ws.Range("I61").Value = 2299.25
$ws.Range("K61").Value = 2299.25
$ws.Range("M61").Value = -2087.25
$ws.Range("H110").Value = 3009.6
$ws.Range("I110").Value = 2887
$ws.Range("K110").Value = 2887
$ws.Range("M110").Value = -842
$ws.Range("H116").Value = 3325.0715
$ws.Range("I116").Value = 3006.625
$ws.Range("J116").Value = 3749.6667
$ws.Range("K116").Value = 3006.625
$ws.Range("L116").Value = 3749.6667
$ws.Range("M116").Value = -712.625
$ws.Range("N116").Value = -8337.6667
$ws.Range("H132").Value = 2778.3064
$ws.Range("I132").Value = 2094.25
$ws.Range("J132").Value = 5123.643
$ws.Range("K132").Value = 6282.75
$ws.Range("L132").Value = 15370.929
$ws.Range("M132").Value = -3752.75
$ws.Range("N132").Value = -20430.929
$ws.Range("H136").Value = 2764.535
$ws.Range("I136").Value = 2299.25
$ws.Range("K136").Value = 6897.75
$ws.Range("M136").Value = -4347.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3325.0715
$ws.Range("I3").Value = 3006.625
$ws.Range("J3").Value = 3749.6667
$ws.Range("K3").Value = 3006.625
$ws.Range("L3").Value = 3749.6667
$ws.Range("M3").Value = -2892.625
$ws.Range("N3").Value = -3977.6667
$ws.Range("H99").Value = 2585.2354
$ws.Range("I99").Value = 1605.5555
$ws.Range("K99").Value = 1605.5555
$ws.Range("M99").Value = -107.5554999999999
$ws.Range("H105").Value = 2519.2703
$ws.Range("I105").Value = 2253.9583
$ws.Range("J105").Value = 3009.077
$ws.Range("K105").Value = 2253.9583
$ws.Range("L105").Value = 3009.077
$ws.Range("M105").Value = -506.9582999999998
$ws.Range("N105").Value = -6503.077
$ws.Range("H134").Value = 1703698.2
$ws.Range("I134").Value = 2166367.8
$ws.Range("J134").Value = 7243.6665
$ws.Range("K134").Value = 6499103.399999999
$ws.Range("L134").Value = 21730.9995
$ws.Range("M134").Value = -6496568.399999999
$ws.Range("N134").Value = -26800.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4531.054
$ws.Range("I31").Value = 1858.0869
$ws.Range("K31").Value = 1858.0869
$ws.Range("M31").Value = -1563.0869
$ws.Range("H34").Value = 4531.054
$ws.Range("I34").Value = 1858.0869
$ws.Range("K34").Value = 1858.0869
$ws.Range("M34").Value = -1656.0869
$ws.Range("H60").Value = 7148.25
$ws.Range("I60").Value = 9546.5
$ws.Range("J60").Value = 4750
$ws.Range("K60").Value = 9546.5
$ws.Range("L60").Value = 4750
$ws.Range("M60").Value = -9035.5
$ws.Range("N60").Value = -5772
$ws.Range("H86").Value = 33588.86
$ws.Range("I86").Value = 23224.363
$ws.Range("J86").Value = 39017.883
$ws.Range("K86").Value = 23224.363
$ws.Range("L86").Value = 39017.883
$ws.Range("M86").Value = -22101.363
$ws.Range("N86").Value = -41263.883
$ws.Range("H89").Value = 33588.86
$ws.Range("I89").Value = 23224.363
$ws.Range("J89").Value = 39017.883
$ws.Range("K89").Value = 116121.815
$ws.Range("L89").Value = 195089.415
$ws.Range("M89").Value = -110505.815
$ws.Range("N89").Value = -206321.415
$ws.Range("H109").Value = 28296
$ws.Range("J109").Value = 28296
$ws.Range("L109").Value = 28296
$ws.Range("N109").Value = -30376
$ws.Range("H122").Value = 2944658.8
$ws.Range("I122").Value = 4548728.5
$ws.Range("J122").Value = 3864.5
$ws.Range("K122").Value = 13646185.5
$ws.Range("L122").Value = 11593.5
$ws.Range("M122").Value = -13643735.5
$ws.Range("N122").Value = -16493.5
$ws.Range("H134").Value = 3862.3076
$ws.Range("I134").Value = 3834.1667
$ws.Range("K134").Value = 11502.5001
$ws.Range("M134").Value = -8967.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3287.1428
$ws.Range("J75").Value = 3168.3333
$ws.Range("L75").Value = 9504.999899999999
$ws.Range("N75").Value = -11500.9999
$ws.Range("H78").Value = 3287.1428
$ws.Range("J78").Value = 3168.3333
$ws.Range("L78").Value = 28514.9997
$ws.Range("N78").Value = -38498.9997
$ws.Range("H97").Value = 539
$ws.Range("H121").Value = 32492.053
$ws.Range("I121").Value = 530
$ws.Range("J121").Value = 36252.293
$ws.Range("K121").Value = 1590
$ws.Range("L121").Value = 108756.879
$ws.Range("M121").Value = -280
$ws.Range("N121").Value = -111376.879
$ws.Range("H131").Value = 1665.1
$ws.Range("I131").Value = 978.5714
$ws.Range("J131").Value = 2034.7693
$ws.Range("K131").Value = 2935.7142
$ws.Range("L131").Value = 6104.3079
$ws.Range("M131").Value = 2104.2858
$ws.Range("N131").Value = -16184.3079
$ws.Range("H132").Value = 3339998.8
$ws.Range("J132").Value = 3339998.8
$ws.Range("L132").Value = 30059989.2
$ws.Range("N132").Value = -30065049.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 24196.75
$ws.Range("I70").Value = 66329.7
$ws.Range("J70").Value = 5045.409
$ws.Range("K70").Value = 66329.7
$ws.Range("L70").Value = 5045.409
$ws.Range("M70").Value = -66059.7
$ws.Range("N70").Value = -5585.409
$ws.Range("H73").Value = 24196.75
$ws.Range("I73").Value = 66329.7
$ws.Range("J73").Value = 5045.409
$ws.Range("K73").Value = 66329.7
$ws.Range("L73").Value = 5045.409
$ws.Range("M73").Value = -65393.7
$ws.Range("N73").Value = -6917.409
$ws.Range("H97").Value = 1078.6471
$ws.Range("I97").Value = 1033.4615
$ws.Range("J97").Value = 1225.5
$ws.Range("K97").Value = 1033.4615
$ws.Range("L97").Value = 1225.5
$ws.Range("M97").Value = -537.4614999999999
$ws.Range("N97").Value = -2217.5
$ws.Range("H126").Value = 3065
$ws.Range("I126").Value = 2847.5
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 8542.5
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -6072.5
$ws.Range("N126").Value = -15440
$ws.Range("H132").Value = 4099.9375
$ws.Range("I132").Value = 3258.5
$ws.Range("K132").Value = 9775.5
$ws.Range("M132").Value = -7245.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 100
$ws.Range("I13").Value = 100
$ws.Range("K13").Value = 100
$ws.Range("M13").Value = 40
$ws.Range("H22").Value = 1833.3077
$ws.Range("I22").Value = 1526
$ws.Range("J22").Value = 2524.75
$ws.Range("K22").Value = 1526
$ws.Range("L22").Value = 2524.75
$ws.Range("M22").Value = -1231
$ws.Range("N22").Value = -3114.75
$ws.Range("H27").Value = 1833.3077
$ws.Range("I27").Value = 1526
$ws.Range("J27").Value = 2524.75
$ws.Range("K27").Value = 1526
$ws.Range("L27").Value = 2524.75
$ws.Range("M27").Value = -1419
$ws.Range("N27").Value = -2738.75
$ws.Range("H136").Value = 2150.2856
$ws.Range("I136").Value = 1925.3334
$ws.Range("K136").Value = 5776.0002
$ws.Range("M136").Value = -3226.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 30045
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 30045
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 30045
$ws.Range("N49").Value = -30505
$ws.Range("M49").ClearContents()
$ws.Range("H55").Value = 17349
$ws.Range("I55").Value = 7047
$ws.Range("J55").Value = 22500
$ws.Range("K55").Value = 7047
$ws.Range("L55").Value = 22500
$ws.Range("M55").Value = -6770
$ws.Range("N55").Value = -23054
$ws.Range("H81").Value = 5106.3
$ws.Range("I81").Value = 4253
$ws.Range("J81").Value = 5959.6
$ws.Range("K81").Value = 8506
$ws.Range("L81").Value = 11919.2
$ws.Range("M81").Value = -7445
$ws.Range("N81").Value = -14041.2
$ws.Range("H84").Value = 5106.3
$ws.Range("I84").Value = 4253
$ws.Range("J84").Value = 5959.6
$ws.Range("K84").Value = 42530
$ws.Range("L84").Value = 59596
$ws.Range("M84").Value = -37226
$ws.Range("N84").Value = -70204
$ws.Range("H107").Value = 398.9
$ws.Range("I107").Value = 325.16666
$ws.Range("K107").Value = 975.4999799999999
$ws.Range("M107").Value = 944.5000200000001
$ws.Range("H132").Value = 1702.4584
$ws.Range("J132").Value = 900
$ws.Range("L132").Value = 2700
$ws.Range("N132").Value = -7760
